$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is being rebuilt from scratch: the old CAS/EC/Name lookup table
# (with heavy borders/fonts for a handful of substances) is replaced by a
# plain CAS-number list (one source worked fully - row 2 still carries the
# EC/Name placeholder cells - the rest only ever populate column A, matching
# "All sources working except for C&L").
$ws.Rows("1:100").Delete()

# Header row: CAS / EC / Name
$ws.Range("A1").Value = "CAS"
$ws.Range("B1").Value = "EC"
$ws.Range("C1").Value = "Name"

# "Name" header cell keeps a left/center aligned, borderless, 11pt black font
$ws.Range("C1").Font.Name = "Aptos Narrow"
$ws.Range("C1").Font.Size = 11
$ws.Range("C1").Font.Color = 0
$ws.Range("C1").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C1").VerticalAlignment = -4108    # xlCenter

# Row 2: the one CAS number that still has matching EC/Name placeholder cells
$ws.Range("A2").Value = "413615-35-7`t"

$ws.Range("A2").Font.Name = "Aptos Narrow"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A2").VerticalAlignment = -4108    # xlCenter
$ws.Range("A2").Borders.Item(7).LineStyle = 1
$ws.Range("A2").Borders.Item(7).Weight = 2
$ws.Range("A2").Borders.Item(10).LineStyle = 1
$ws.Range("A2").Borders.Item(10).Weight = 2
$ws.Range("A2").Borders.Item(8).LineStyle = 1
$ws.Range("A2").Borders.Item(8).Weight = 2
$ws.Range("A2").Borders.Item(9).LineStyle = 1
$ws.Range("A2").Borders.Item(9).Weight = 2

$ws.Range("B2").Font.Name = "Aptos Narrow"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B2").VerticalAlignment = -4108    # xlCenter
$ws.Range("B2").Borders.Item(7).LineStyle = 1
$ws.Range("B2").Borders.Item(7).Weight = 2
$ws.Range("B2").Borders.Item(10).LineStyle = 1
$ws.Range("B2").Borders.Item(10).Weight = -4138
$ws.Range("B2").Borders.Item(8).LineStyle = 1
$ws.Range("B2").Borders.Item(8).Weight = 2
$ws.Range("B2").Borders.Item(9).LineStyle = 1
$ws.Range("B2").Borders.Item(9).Weight = 2

$ws.Range("C2").Font.Name = "Aptos Narrow"
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.Color = 0
$ws.Range("C2").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C2").VerticalAlignment = -4108    # xlCenter
$ws.Range("C2").Borders.Item(7).LineStyle = 1
$ws.Range("C2").Borders.Item(7).Weight = -4138
$ws.Range("C2").Borders.Item(10).LineStyle = 1
$ws.Range("C2").Borders.Item(10).Weight = 2
$ws.Range("C2").Borders.Item(8).LineStyle = 1
$ws.Range("C2").Borders.Item(8).Weight = 2
$ws.Range("C2").Borders.Item(9).LineStyle = 1
$ws.Range("C2").Borders.Item(9).Weight = 2

# Remaining rows: CAS numbers only, plain cells in column A
$ws.Range("A3").Value = "1166-46-7"
$ws.Range("A4").Value = "141-62-8"
$ws.Range("A5").Value = "141-63-9"
$ws.Range("A6").Value = "85877-79-8"
$ws.Range("A7").Value = "75-09-2"
$ws.Range("A8").Value = "597-82-0"
$ws.Range("A9").Value = "72963-72-5"
$ws.Range("A10").Value = "1195-32-0"
$ws.Range("A11").Value = "215-724-4"

# A9 keeps a left/center aligned cell with a partial box border
$ws.Range("A9").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A9").VerticalAlignment = -4108    # xlCenter
$ws.Range("A9").Borders.Item(7).LineStyle = 1
$ws.Range("A9").Borders.Item(7).Weight = 2
$ws.Range("A9").Borders.Item(10).LineStyle = 1
$ws.Range("A9").Borders.Item(10).Weight = -4138
$ws.Range("A9").Borders.Item(8).LineStyle = 1
$ws.Range("A9").Borders.Item(8).Weight = 2
$ws.Range("A9").Borders.Item(9).LineStyle = 1
$ws.Range("A9").Borders.Item(9).Weight = 2

# Matches the saved selection in the target file
$ws.Range("A11").Select()
